# Weekly update: two new daily observations are inserted at the top of the
# "Cilantro" price series (before the old row 88), pushing the remainder of
# the series down by two rows. The newly inserted rows start out as copies
# of the data that lands right below them (the historical next-oldest
# observations) and then get their own date/volume corrected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 88 - everything from
# row 88 downward (including the former last row, 254) shifts down by two,
# which is exactly what the diff shows happening to rows 90-256.
$ws.Rows.Item(88).Insert()
$ws.Rows.Item(89).Insert()

# Seed the two new rows with the content now sitting just below them (the
# rows that used to be 88 and 89 before the insert), then correct the
# date (column D) and volume (column J) for the new week.
$ws.Rows.Item(90).Copy()
$ws.Rows.Item(88).PasteSpecial()

$ws.Rows.Item(91).Copy()
$ws.Rows.Item(89).PasteSpecial()

$ws.Cells.Item(88, 4).Value = 44915
$ws.Cells.Item(88, 10).Value = 2000

$ws.Cells.Item(89, 4).Value = 44915
$ws.Cells.Item(89, 10).Value = 1000
